$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.042494666666667
$ws.Range("H2").Value = 3.127484
$ws.Range("I2").Value = 0.0007670466909205676
$ws.Range("J2").Value = 0.0007670466909205677
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.1152498703893333
$ws.Range("R2").Value = 1.037248833504
$ws.Range("S2").Value = 0.000008638859994717835
$ws.Range("T2").Value = 0.000008638859994717838

$ws.Range("G3").Value = 1.042494666666667
$ws.Range("H3").Value = 3.127484
$ws.Range("I3").Value = 0.0007670466909205676
$ws.Range("J3").Value = 0.0007670466909205677
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 9.395035605623109
$ws.Range("R3").Value = 84.555320450608
$ws.Range("S3").Value = 0.0007042298352977492
$ws.Range("T3").Value = 0.0007042298352977493

$ws.Range("G4").Value = 1.042494666666667
$ws.Range("H4").Value = 3.127484
$ws.Range("I4").Value = 0.0007670466909205676
$ws.Range("J4").Value = 0.0007670466909205677
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 0.7227813597986665
$ws.Range("R4").Value = 6.505032238187999
$ws.Range("S4").Value = 0.00005417799562810057
$ws.Range("T4").Value = 0.00005417799562810059

$ws.Range("I5").Value = 0.9658609009611662
$ws.Range("J5").Value = 0.9658609009611662
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 145.1219918781013
$ws.Range("R5").Value = 1306.097926902912
$ws.Range("S5").Value = 0.01087800416394679
$ws.Range("T5").Value = 0.0108780041639468

$ws.Range("I6").Value = 0.9658609009611662
$ws.Range("J6").Value = 0.9658609009611662
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("S6").Value = 0.8867622678719767
$ws.Range("T6").Value = 0.8867622678719768

$ws.Range("I7").Value = 0.9658609009611662
$ws.Range("J7").Value = 0.9658609009611662
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 910.1222437127625
$ws.Range("R7").Value = 8191.100193414863
$ws.Range("S7").Value = 0.06822062892524267
$ws.Range("T7").Value = 0.06822062892524268

$ws.Range("G8").Value = 45.356022
$ws.Range("H8").Value = 136.068066
$ws.Range("I8").Value = 0.03337205234791334
$ws.Range("J8").Value = 0.03337205234791334
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 5.014198944144
$ws.Range("R8").Value = 45.127790497296
$ws.Range("S8").Value = 0.0003758525933069605
$ws.Range("T8").Value = 0.0003758525933069606

$ws.Range("G9").Value = 45.356022
$ws.Range("H9").Value = 136.068066
$ws.Range("I9").Value = 0.03337205234791334
$ws.Range("J9").Value = 0.03337205234791334
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 408.751675422888
$ws.Range("R9").Value = 3678.765078805991
$ws.Range("S9").Value = 0.03063906696515898
$ws.Range("T9").Value = 0.03063906696515898

$ws.Range("G10").Value = 45.356022
$ws.Range("H10").Value = 136.068066
$ws.Range("I10").Value = 0.03337205234791334
$ws.Range("J10").Value = 0.03337205234791334
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 31.44619181701799
$ws.Range("R10").Value = 283.0157263531619
$ws.Range("S10").Value = 0.002357132789447396
$ws.Range("T10").Value = 0.002357132789447397
